$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# New header cells H1:L1 need the same header formatting (bold, bordered,
# centered/top-aligned) already used by the existing header cells, so
# copy the format from an existing header cell before setting values.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("H1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F1").Value = "metric"
$ws.Range("G1").Value = "period"
$ws.Range("H1").Value = "tool_value"
$ws.Range("I1").Value = "target_value"
$ws.Range("J1").Value = "comparator"
$ws.Range("K1").Value = "result"
$ws.Range("L1").Value = "reason"

# --- Row 2 ---
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = "Success"
$ws.Range("L2").Value = "Domain found in contacts"

# --- Row 3 ---
$ws.Range("C3").Value = "mail_tool"
# D3/E3 were already empty and remain untouched by the diff
$ws.Range("F3").Value = $null
$ws.Range("G3").Value = $null
$ws.Range("H3").Value = $null
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = $null
$ws.Range("K3").Value = "Failed"
$ws.Range("L3").Value = "Invalid or missing email format"

# --- Row 4 ---
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = "Success"
$ws.Range("L4").Value = "Domain found in contacts"

# --- Row 5 ---
$ws.Range("C5").Value = "mail_tool"
# D5/E5 were already empty and remain untouched by the diff
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = $null
$ws.Range("H5").Value = $null
$ws.Range("I5").Value = $null
$ws.Range("J5").Value = $null
$ws.Range("K5").Value = "Failed"
$ws.Range("L5").Value = "Invalid or missing email format"

# --- Row 6 (new) ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Aggregate monthly data for March 2025 should be less than 1200"
$ws.Range("C6").Value = "monthly_tool"
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = $null
$ws.Range("F6").Value = $null
$ws.Range("G6").Value = "month:2025-03"
$ws.Range("H6").Value = 6190
$ws.Range("I6").Value = 2025
$ws.Range("J6").Value = "lt"
$ws.Range("K6").Value = "Failed"
$ws.Range("L6").Value = "tool_value=6190.0, target_value=2025.0, comparator=lt, period=month:2025-03, metric=None"

# --- Row 7 (new) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "CSR supply in April 2025 equals 450"
$ws.Range("C7").Value = "monthly_tool"
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = "csr_supply"
$ws.Range("G7").Value = "month:2025-04"
$ws.Range("H7").Value = 5850
$ws.Range("I7").Value = 2025
$ws.Range("J7").Value = "eq"
$ws.Range("K7").Value = "Failed"
$ws.Range("L7").Value = "tool_value=5850.0, target_value=2025.0, comparator=eq, period=month:2025-04, metric=csr_supply"

# --- Row 8 (new) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Total for 2025-03 must be >= 1000"
$ws.Range("C8").Value = "monthly_tool"
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = "month:2025-03"
$ws.Range("H8").Value = 6190
$ws.Range("I8").Value = 2025
$ws.Range("J8").Value = "ge"
$ws.Range("K8").Value = "Success"
$ws.Range("L8").Value = "tool_value=6190.0, target_value=2025.0, comparator=ge, period=month:2025-03, metric=None"

# --- Row 9 (new) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Spend on this item was this week 932 dollars"
$ws.Range("C9").Value = "monthly_tool"
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = "spend"
$ws.Range("G9").Value = "week:2025-93"
$ws.Range("H9").Value = $null
$ws.Range("I9").Value = 932
$ws.Range("J9").Value = $null
$ws.Range("K9").Value = "Failed"
$ws.Range("L9").Value = "No data"
